$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cheques")

# Row 4
$ws.Range("A4").Value = 731950
$ws.Range("B4").Value = "'7556522"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "'4"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "EDILIO DEL TRANSITO ALVAREZ TALAMILLA"
$ws.Range("E4").Value = "'29743125"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "'96509660"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "'210021989816"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = "'146687"
$ws.Range("H4").Style = "Normal"
$ws.Range("I4").Value = "'2024-09-27"
$ws.Range("I4").Style = "Normal"
$ws.Range("J4").Value = "AUTOMOTRIZ"
$ws.Range("K4").Value = "JUDICIAL"
$ws.Range("L4").Value = "1 CUOTA"
$ws.Range("M4").Style = "Normal"

# Row 5
$ws.Range("A5").Value = 731951
$ws.Range("B5").Value = "'14191212"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = "'7"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "PATRICIO ALEJANDRO HURTADO ALVAREZ"
$ws.Range("E5").Value = "'29743125"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "'96509660"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "'230017928138"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = "'165624"
$ws.Range("H5").Style = "Normal"
$ws.Range("I5").Value = "'2024-05-09"
$ws.Range("I5").Style = "Normal"
$ws.Range("J5").Value = "RENEGOCIACION"
$ws.Range("K5").Value = "VIGENTE"
$ws.Range("L5").Value = "1 CUOTA"
$ws.Range("M5").Style = "Normal"

# Row 6
$ws.Range("A6").Value = 731952
$ws.Range("B6").Value = "'11671345"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = "'4"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "FIERRO REYES CARLOS YASHIN"
$ws.Range("E6").Value = "'29743125"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "'96509660"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "'206005354319"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = "'281884"
$ws.Range("H6").Style = "Normal"
$ws.Range("I6").Value = "'2024-06-05"
$ws.Range("I6").Style = "Normal"
$ws.Range("J6").Value = "CONSUMO"
$ws.Range("K6").Value = "VIGENTE"
$ws.Range("L6").Value = "1 CUOTA"
$ws.Range("M6").Style = "Normal"

# Row 7
$ws.Range("A7").Value = 731953
$ws.Range("B7").Value = "'17199914"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").Value = "'6"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = "CASTRO ARRIAGADA EVELIN NATALI"
$ws.Range("E7").Value = "'61682420"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "'96509669"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Value = "'558330389530"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value = "'66879"
$ws.Range("H7").Style = "Normal"
$ws.Range("I7").Value = "'2024-10-21"
$ws.Range("I7").Style = "Normal"
$ws.Range("J7").Style = "Normal"
$ws.Range("K7").Value = "CASTIGO"
$ws.Range("L7").Value = "1 CUOTA"
$ws.Range("M7").Style = "Normal"

# J2 empty placeholder cell
$ws.Range("J2").Style = "Normal"
